# adding DaemonSet and Logging
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the three new rows with kind / apiVersion pairs
$ws.Range("A8").Value = "LimitRange"
$ws.Range("B8").Value = "v1"

$ws.Range("A9").Value = "ResourceQuota"
$ws.Range("B9").Value = "v1"

$ws.Range("A10").Value = "DaemonSet"
$ws.Range("B10").Value = "apps/v1"

# Update the active selection as recorded in the saved file
$ws.Range("J2").Select()
